{"js": "// Update the date paragraph and the division-fact table cells to the new\n// content, preserving all existing formatting (fonts, sizes, alignment).\n\n// 1) Date heading paragraph: \"2025-06-30 Monday\" -> \"2025-07-01 Tuesday\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  paragraphs.items[0].insertText(\"2025-07-01 Tuesday\", \"Replace\");\n}\n\n// 2) The worksheet table: 5 populated rows (at table-row indices 0, 4, 8,\n//    12, 16 - interleaved with blank spacer rows) x 5 columns of\n//    \"A\u00f7B=C, D\" style division facts. Replace each populated cell's text\n//    with its new value, addressed by (row, col) so the duplicate old\n//    value \"59\u00f74=14, 3\" (rows 4 and 8, col 1 vs col 0) still resolves to\n//    the correct distinct replacement.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// newValues[tableRowIndex] = [col0, col1, col2, col3, col4]\nconst newValues = {\n  0: [\"66\u00f76=11, 0\", \"76\u00f79=8, 4\", \"27\u00f76=4, 3\", \"19\u00f72=9, 1\", \"61\u00f75=12, 1\"],\n  4: [\"70\u00f73=23, 1\", \"84\u00f72=42, 0\", \"82\u00f77=11, 5\", \"36\u00f73=12, 0\", \"32\u00f78=4, 0\"],\n  8: [\"73\u00f76=12, 1\", \"80\u00f76=13, 2\", \"92\u00f74=23, 0\", \"51\u00f75=10, 1\", \"90\u00f77=12, 6\"],\n  12: [\"32\u00f73=10, 2\", \"15\u00f74=3, 3\", \"64\u00f78=8, 0\", \"69\u00f75=13, 4\", \"29\u00f74=7, 1\"],\n  16: [\"25\u00f74=6, 1\", \"47\u00f79=5, 2\", \"70\u00f76=11, 4\", \"33\u00f79=3, 6\", \"71\u00f74=17, 3\"]\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = Number(rowIndex);\n  const rowVals = newValues[rowIndex];\n  for (let c = 0; c < rowVals.length; c++) {\n    table.getCell(r, c).value = rowVals[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the division-fact table cells to the new\n# content, preserving all existing formatting (fonts, sizes, alignment).\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph: \"2025-06-30 Monday\" -> \"2025-07-01 Tuesday\".\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-01 Tuesday\"\n\n# 2) The worksheet table: 5 populated rows (1-based table rows 1, 5, 9, 13,\n#    17 - interleaved with blank spacer rows) x 5 columns of \"A\u00f7B=C, D\"\n#    style division facts. Replace each populated cell's text with its new\n#    value, addressed by (row, col) so the duplicate old value\n#    \"59\u00f74=14, 3\" (rows 5 and 9, col 2 vs col 1) still resolves to the\n#    correct distinct replacement.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"66\u00f76=11, 0\", \"76\u00f79=8, 4\", \"27\u00f76=4, 3\", \"19\u00f72=9, 1\", \"61\u00f75=12, 1\")\n    5  = @(\"70\u00f73=23, 1\", \"84\u00f72=42, 0\", \"82\u00f77=11, 5\", \"36\u00f73=12, 0\", \"32\u00f78=4, 0\")\n    9  = @(\"73\u00f76=12, 1\", \"80\u00f76=13, 2\", \"92\u00f74=23, 0\", \"51\u00f75=10, 1\", \"90\u00f77=12, 6\")\n    13 = @(\"32\u00f73=10, 2\", \"15\u00f74=3, 3\", \"64\u00f78=8, 0\", \"69\u00f75=13, 4\", \"29\u00f74=7, 1\")\n    17 = @(\"25\u00f74=6, 1\", \"47\u00f79=5, 2\", \"70\u00f76=11, 4\", \"33\u00f79=3, 6\", \"71\u00f74=17, 3\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowVals = $newValues[$rowIndex]\n    for ($c = 0; $c -lt $rowVals.Length; $c++) {\n        $t.Cell($rowIndex, $c + 1).Range.Text = $rowVals[$c]\n    }\n}\n"}
